$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New user stories / tasks appended below the existing backlog (rows 56-85).
$ws.Range("A56").Value = "S7: As a user I want to enhance login functinality so that user can login as trader or broker role"
$ws.Range("B56").Value = "S7T1: As a developer I need to add radio buttons for trader and broker so that user can able to select one."
$ws.Range("B57").Value = "S7T2: As a developer I need to "
$ws.Range("A66").Value = "S8: As a user I want  accept order functinality so that I can accept the executed order."
$ws.Range("A76").Value = "S9: As a user I want execute order functionality so that user can execute order "
$ws.Range("A85").Value = "S10: As a user I want logout functinality."

# Match the author's final viewport/selection state.
$ws.Range("B57").Select()
$excel.ActiveWindow.ScrollRow = 41
